$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")

# Insert a new row at sheet row 212 (pushes the old row 212.. down by one,
# and the ListObject's underlying sheet dimension grows automatically).
$ws.Rows.Item(212).Insert()

# The inserted row picked up default (un-bordered) formatting; copy the
# cell formatting from row 214 (a plain data row that already has the
# styles we want: date style in A, box style in B/C/D/E/F/G/H/I/J/K) down
# onto the blank new row 212.
$ws.Range("A214:K214").Copy()
$ws.Range("A212:K212").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new leave-card entry. (Set the REMARKS string before the
# PARTICULARS string so new shared-string entries are appended in the
# same order the original author typed them: period text, then code.)
$ws.Range("K212").Value = "01/02-12/2024"
$ws.Range("B212").Value = "VL(9-0-0)"
$ws.Range("D212").Value = 9
$ws.Range("G212").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# Row 211 (CTO) now also has an EARNED value recorded.
$ws.Range("C211").Value = 1.25

# Grow the table boundary to include the newly appended row at the end
# (shifted down from the original last row).
$lo.Resize($ws.Range("A8:K345"))

# The row that fell off the bottom of the table during the insert/resize
# picked up a calculated-column formula that was (momentarily) out of
# table bounds and evaluated to an error; re-stamp it now that the row is
# back inside Table1 so it recomputes to the normal blank string.
$ws.Range("G345").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# Restore the active selection to match the edited cell.
$ws.Range("C214").Select()
